# Update the data value in A1 from 2 to 5, and move the active selection
# from D4 to D5 (matches commit: "Обновил данные на 5").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5
$ws.Range("D5").Select()
